$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.316.56"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "3.034.00"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("E4").Value = "  +0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "579.02"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.54%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "167.98"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.78%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "3.031.23"
$ws.Range("E8").Value = "  +1.01%  "
$ws.Range("E9").Value = "  +0.51%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "6.66"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("E11").Value = "  -1.62%  "
$ws.Range("E12").Value = "  +7.37%  "
$ws.Range("E13").Value = "  -1.69%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "36.65"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +5.99%  "
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "66.317.89"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").Value = "3.536.81"
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("E18").Value = "  +4.53%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "16.56"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +19.48%  "
$ws.Range("D20").Value = "3.032.82"
$ws.Range("E20").Value = "  +0.88%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "466.38"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +2.90%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.712"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +3.89%  "
$ws.Range("E23").Value = "  +0.81%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "83.07"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.90%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "12.72"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +3.73%  "
$ws.Range("E26").Value = "  -1.22%  "
$ws.Range("E27").Value = "  -2.54%  "
$ws.Range("E28").Value = "  -0.01%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "8.22"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.00%  "
$ws.Range("E30").Value = "  +0.63%  "
$ws.Range("E31").Value = "  +1.19%  "
$ws.Range("E32").Value = "  +6.44%  "
$ws.Range("D33").Value = "0.0₃0992"
$ws.Range("E33").Value = "  -4.60%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "28.19"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +3.57%  "
$ws.Range("E35").Value = "  +0.09%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.993"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("E37").Value = "  +0.90%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "48.48"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +10.04%  "
$ws.Range("E39").Value = "  -0.68%  "
$ws.Range("E40").Value = "  +2.32%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("E42").Value = "  -0.72%  "
$ws.Range("E43").Value = "  +2.43%  "
$ws.Range("E44").Value = "  -3.95%  "
$ws.Range("E45").Value = "  +0.33%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "380.29"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -5.56%  "
$ws.Range("D47").Value = "2.716.41"
$ws.Range("E47").Value = "  -1.90%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "133.76"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("E49").Value = "  +0.01%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "24.50"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +2.86%  "
$ws.Range("E51").Value = "  +3.99%  "
